$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ordering / content of the Fuel -> Colour lookup table.
# Rows 2-19 lose their border-style formatting (ClearFormats -> no "s" attr),
# rows 20-23 keep their existing style.
$rows = @(
    @{ Row = 2;  Fuel = "Electricity";         Color = "#0086be"; Styled = $false },
    @{ Row = 3;  Fuel = "Petrol";               Color = "#a7a9ac"; Styled = $false },
    @{ Row = 4;  Fuel = "Diesel";               Color = "#58595b"; Styled = $false },
    @{ Row = 5;  Fuel = "Hydrogen";             Color = "#8a73b4"; Styled = $false },
    @{ Row = 6;  Fuel = "Fuel Oil";             Color = "#974a21"; Styled = $false },
    @{ Row = 7;  Fuel = "Coal";                 Color = "#231f20"; Styled = $false },
    @{ Row = 8;  Fuel = "Natural Gas";          Color = "#f0575b"; Styled = $false },
    @{ Row = 9;  Fuel = "Wood";                 Color = "#53817a"; Styled = $false },
    @{ Row = 10; Fuel = "Geothermal";           Color = "#4e2e8e"; Styled = $false },
    @{ Row = 11; Fuel = "LPG";                  Color = "#c12025"; Styled = $false },
    @{ Row = 12; Fuel = "Biodiesel";            Color = "#00435b"; Styled = $false },
    @{ Row = 13; Fuel = "Drop-in Diesel";       Color = "#00435b"; Styled = $false },
    @{ Row = 14; Fuel = "Biogas";               Color = "#00435b"; Styled = $false },
    @{ Row = 15; Fuel = "Hydro";                Color = "#49bee5"; Styled = $false },
    @{ Row = 16; Fuel = "Solar";                Color = "#00af8c"; Styled = $false },
    @{ Row = 17; Fuel = "Waste Incineration";   Color = "#6b0d0e"; Styled = $false },
    @{ Row = 18; Fuel = "Wind";                 Color = "#ffc808"; Styled = $false },
    @{ Row = 19; Fuel = "Jet Fuel";             Color = "#f57e20"; Styled = $false },
    @{ Row = 20; Fuel = "Tui";                  Color = "#164057"; Styled = $true },
    @{ Row = 21; Fuel = "Kea";                  Color = "#00af8c"; Styled = $true },
    @{ Row = 22; Fuel = "Fossil Fuel";          Color = "#414042"; Styled = $true },
    @{ Row = 23; Fuel = "Renewable";            Color = "#00af8c"; Styled = $true }
)

foreach ($item in $rows) {
    $r = $item.Row
    $cellA = $ws.Cells.Item($r, 1)
    $cellB = $ws.Cells.Item($r, 2)
    $cellA.Value = $item.Fuel
    $cellB.Value = $item.Color
    if (-not $item.Styled) {
        $cellA.ClearFormats()
        $cellB.ClearFormats()
        $cellA.Value = $item.Fuel
        $cellB.Value = $item.Color
    }
}

$null = $ws.Range("I5").Select()
